$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 18222.416
$ws.Range("I88").Value = 1751.5
$ws.Range("J88").Value = 21516.6
$ws.Range("K88").Value = 1751.5
$ws.Range("L88").Value = 21516.6
$ws.Range("M88").Value = -1345.5
$ws.Range("N88").Value = -22328.6

$ws.Range("H91").Value = 18222.416
$ws.Range("I91").Value = 1751.5
$ws.Range("J91").Value = 21516.6
$ws.Range("K91").Value = 1751.5
$ws.Range("L91").Value = 21516.6
$ws.Range("M91").Value = -347.5
$ws.Range("N91").Value = -24324.6

$ws.Range("H92").Value = 1518.4706
$ws.Range("I92").Value = 728.0909
$ws.Range("J92").Value = 2967.5
$ws.Range("K92").Value = 728.0909
$ws.Range("L92").Value = 2967.5
$ws.Range("M92").Value = 519.9091
$ws.Range("N92").Value = -5463.5

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 1000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -1232

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H9").Value = 10006
$ws.Range("J9").Value = 10006
$ws.Range("L9").Value = 10006
$ws.Range("N9").Value = -10346

$ws.Range("H20").Value = 10006
$ws.Range("J20").Value = 10006
$ws.Range("L20").Value = 10006
$ws.Range("N20").Value = -10546

$ws.Range("H23").Value = 20000
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20518

$ws.Range("H37").Value = 40000
$ws.Range("J37").Value = 40000
$ws.Range("L37").Value = 40000
$ws.Range("N37").Value = -40546

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H55").Value = 3000
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -2685
$ws.Range("N55").ClearContents()

$ws.Range("H63").Value = 3131.6667
$ws.Range("I63").Value = 2898.125
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2898.125
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -2212.125
$ws.Range("N63").Value = -6372

$ws.Range("H66").Value = 3131.6667
$ws.Range("I66").Value = 2898.125
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 14490.625
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -11058.625
$ws.Range("N66").Value = -31864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 765.7
$ws.Range("I22").Value = 533.8
$ws.Range("J22").Value = 997.6
$ws.Range("K22").Value = 533.8
$ws.Range("L22").Value = 997.6
$ws.Range("M22").Value = -360.8
$ws.Range("N22").Value = -1343.6

$ws.Range("H94").Value = 842.2895
$ws.Range("I94").Value = 733.8276
$ws.Range("J94").Value = 1191.7778
$ws.Range("K94").Value = 733.8276
$ws.Range("L94").Value = 1191.7778
$ws.Range("M94").Value = -282.8276
$ws.Range("N94").Value = -2093.7778

$ws.Range("H132").Value = 44510
$ws.Range("J132").Value = 44510
$ws.Range("L132").Value = 44510
$ws.Range("N132").Value = -54630

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 149
$ws.Range("I7").Value = 44.75
$ws.Range("J7").Value = 357.5
$ws.Range("K7").Value = 44.75
$ws.Range("L7").Value = 357.5
$ws.Range("M7").Value = 68.25
$ws.Range("N7").Value = -583.5

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("N11").ClearContents()

$ws.Range("H122").Value = 1226.2667
$ws.Range("I122").Value = 1001.56757
$ws.Range("J122").Value = 2265.5
$ws.Range("K122").Value = 3004.70271
$ws.Range("L122").Value = 6796.5
$ws.Range("M122").Value = -554.70271
$ws.Range("N122").Value = -11696.5

$ws.Range("H132").Value = 74927
$ws.Range("I132").Value = 3670
$ws.Range("J132").Value = 203189.6
$ws.Range("K132").Value = 11010
$ws.Range("L132").Value = 609568.8
$ws.Range("M132").Value = -8480
$ws.Range("N132").Value = -614628.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 5161987.5
$ws.Range("I11").Value = 6285784
$ws.Range("J11").Value = 666800
$ws.Range("K11").Value = 18857352
$ws.Range("L11").Value = 2000400
$ws.Range("M11").Value = -18857212
$ws.Range("N11").Value = -2000680

$ws.Range("H17").Value = 885.7143
$ws.Range("I17").Value = 600
$ws.Range("K17").Value = 1800
$ws.Range("M17").Value = -1631

$ws.Range("H34").Value = 1066.6666
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 1228.5714
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 3685.7142
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -3853.7142

$ws.Range("H37").Value = 200000
$ws.Range("J37").Value = 200000
$ws.Range("L37").Value = 600000
$ws.Range("N37").Value = -600224

$ws.Range("H39").Value = 2000
$ws.Range("J39").Value = 2000
$ws.Range("L39").Value = 6000
$ws.Range("N39").Value = -6588

$ws.Range("H40").Value = 273.42856
$ws.Range("I40").Value = 40
$ws.Range("J40").Value = 584.6667
$ws.Range("K40").Value = 160
$ws.Range("L40").Value = 2338.6668
$ws.Range("M40").Value = -91
$ws.Range("N40").Value = -2476.6668

$ws.Range("H55").Value = 3135.0588
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 3206
$ws.Range("K55").Value = 6000
$ws.Range("L55").Value = 9618
$ws.Range("M55").Value = -5823
$ws.Range("N55").Value = -9972

$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()

$ws.Range("H133").Value = 4329.4
$ws.Range("I133").Value = 3786
$ws.Range("J133").Value = 4872.8
$ws.Range("K133").Value = 11358
$ws.Range("L133").Value = 14618.4
$ws.Range("M133").Value = -6298
$ws.Range("N133").Value = -24738.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2041.1
$ws.Range("I97").Value = 2287.2856
$ws.Range("J97").Value = 1466.6666
$ws.Range("K97").Value = 2287.2856
$ws.Range("L97").Value = 1466.6666
$ws.Range("M97").Value = -1791.2856
$ws.Range("N97").Value = -2458.6666

$ws.Range("H132").Value = 54596.258
$ws.Range("I132").Value = 41521.73
$ws.Range("J132").Value = 80745.30499999999
$ws.Range("K132").Value = 124565.19
$ws.Range("L132").Value = 242235.915
$ws.Range("M132").Value = -122035.19
$ws.Range("N132").Value = -247295.915

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 821.36365
$ws.Range("I22").Value = 436.375
$ws.Range("J22").Value = 1041.3572
$ws.Range("K22").Value = 436.375
$ws.Range("L22").Value = 1041.3572
$ws.Range("M22").Value = -141.375
$ws.Range("N22").Value = -1631.3572

$ws.Range("H27").Value = 821.36365
$ws.Range("I27").Value = 436.375
$ws.Range("J27").Value = 1041.3572
$ws.Range("K27").Value = 436.375
$ws.Range("L27").Value = 1041.3572
$ws.Range("M27").Value = -329.375
$ws.Range("N27").Value = -1255.3572

$ws.Range("H132").Value = 171167.17
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 171167.17
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 513501.51
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -518561.51
